$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the final meeting row (row 9) attendance for both members.
$ws.Range("B9").Value = "Y"
$ws.Range("C9").Value = "Y"

# Update the saved view/selection state to match the latest edit position.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("C9").Select()
